# Add group-management / support-chat / settings phrases to the "test"
# localization sheet (column A), replacing the former last row
# ("Настройки ") with a longer run of new bot strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test")

$ws.Range("A42").Value = "Изменение настроек"
$ws.Range("A43").Value = "Список групп"
$ws.Range("A44").Value = "Добавить группу"
$ws.Range("A45").Value = "Добавте меня в группу и используйте команду /chat_id чтобы я связался с группой"
$ws.Range("A46").Value = "Укажите задержку с которойбудут приходить сообщения от клиентов в группу"
$ws.Range("A47").Value = "Укажите тематику группы"
$ws.Range("A48").Value = "Ошибка, попробуйте повторно"
$ws.Range("A49").Value = "Группа создана"
$ws.Range("A50").Value = "Укажите новую задержку"
$ws.Range("A51").Value = "Изменения внесены успешно"
$ws.Range("A52").Value = "Укажите новую тему"
$ws.Range("A53").Value = "Точно удалить группу?"
$ws.Range("A54").Value = "Удаление прошло успешно"
$ws.Range("A55").Value = "Выберите по какой тематике искать юриста"
$ws.Range("A56").Value = "Спасибо, ваша заявка отправлена в группы."
# Row 56 did not exist before (gap in the sheet); give it the same
# wrap-text cell style ("s=1") used by every other row in this column.
$ws.Range("A56").WrapText = $true
$ws.Range("A57").Value = "Опишите свою проблему"
$ws.Range("A58").Value = "Вот список ваших запровсов"
$ws.Range("A59").Value = "Текст поиска:"
$ws.Range("A60").Value = "Удалить поиск"
$ws.Range("A61").Value = "Поиск по запросу удален"
$ws.Range("A62").Value = "Введите код клиента"
$ws.Range("A63").Value = "Такого кода нет в нашей базе"
$ws.Range("A64").Value = "Пожалуйста введите первое сообщение для клиента "
$ws.Range("A65").Value = "Пользователь №"
$ws.Range("A66").Value = "Отправить сообщение"
$ws.Range("A67").Value = "Закончить разговор"
$ws.Range("A68").Value = "Разговор окончен"
$ws.Range("A69").Value = "Выберите категорию настроек"
$ws.Range("A70").Value = "Изменение текстов бота"
$ws.Range("A71").Value = "Изменение пароля"
$ws.Range("A72").Value = "Пожалуйста введите новый пароль"
$ws.Range("A73").Value = "Изменение прошло успешно"
$ws.Range("A74").Value = "Вот старые фразы, измените их и отошлите мне фаил вновь."

# The two freshly-added, long wrapped-text rows render taller in Excel.
$ws.Range("A45").RowHeight = 30
$ws.Range("A46").RowHeight = 30

# Scroll the saved view down to where the new rows are, keep the same
# active cell selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 37
$ws.Range("A42").Select()

Write-Output "done"
